# Commit: Mon, Apr 06, 2020 11:05:53 PM
#
# Change the table style ("Table Design" gallery pick) on the single
# table in the deck (slide 16) from the custom "Table_0" style
# ({D8C75BBE-99FA-4111-8F3E-2FE490FA3C6C}) to the built-in style
# {A8DA4D68-D55C-4764-BAB5-A33261A8120E} (PowerPoint's "Medium Style 2 -
# Accent 1"), exactly as if the user had clicked a new style thumbnail
# in the Table Design ribbon tab.

$p = $ppt.ActivePresentation

$targetStyleId  = "{A8DA4D68-D55C-4764-BAB5-A33261A8120E}"
$previousStyleId = "{D8C75BBE-99FA-4111-8F3E-2FE490FA3C6C}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            if ($shape.Table.Style -eq $previousStyleId) {
                $shape.Table.ApplyStyle($targetStyleId)
            }
        }
    }
}
